$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'319.96"
$ws.Range("E2").Value = "'2.36%"
$ws.Range("G2").Value = "'2"
$ws.Range("D3").Value = "'39.84"
$ws.Range("E3").Value = "'5.74%"
$ws.Range("G3").Value = "'2"
$ws.Range("D4").Value = "'5.228"
$ws.Range("E4").Value = "'1.85%"
$ws.Range("G4").Value = "'2"
$ws.Range("D5").Value = "'0.08140"
$ws.Range("E5").Value = "'3.09%"
$ws.Range("G5").Value = "'2"
$ws.Range("D6").Value = "'8.601"
$ws.Range("E6").Value = "'4.01%"
$ws.Range("G6").Value = "'2"
$ws.Range("D7").Value = "'1.907"
$ws.Range("E7").Value = "'0.57%"
$ws.Range("G7").Value = "'2"
$ws.Range("D8").Value = "'2.999"
$ws.Range("E8").Value = "'0.12%"
$ws.Range("G8").Value = "'2"
$ws.Range("D9").Value = "'0.9397"
$ws.Range("E9").Value = "'2.06%"
$ws.Range("G9").Value = "'2"
$ws.Range("D10").Value = "'0.1276"
$ws.Range("E10").Value = "'9.62%"
$ws.Range("G10").Value = "'2"
$ws.Range("D11").Value = "'0.1961"
$ws.Range("E11").Value = "'2.97%"
$ws.Range("G11").Value = "'2"
$ws.Range("D12").Value = "'0.09166"
$ws.Range("E12").Value = "'0.58%"
$ws.Range("G12").Value = "'2"
$ws.Range("D13").Value = "'0.03379"
$ws.Range("E13").Value = "'2.04%"
$ws.Range("G13").Value = "'2"
$ws.Range("D14").Value = "'0.09514"
$ws.Range("E14").Value = "'-0.98%"
$ws.Range("G14").Value = "'2"
$ws.Range("D15").Value = "'0.001397"
$ws.Range("E15").Value = "'1.18%"
$ws.Range("G15").Value = "'2"
$ws.Range("D16").Value = "'0.005927"
$ws.Range("E16").Value = "'1.01%"
$ws.Range("G16").Value = "'2"
$ws.Range("D17").Value = "'3.366"
$ws.Range("E17").Value = "'-5.30%"
$ws.Range("G17").Value = "'2"
$ws.Range("D18").Value = "'4.505"
$ws.Range("E18").Value = "'2.00%"
$ws.Range("G18").Value = "'2"
$ws.Range("D19").Value = "'0.3536"
$ws.Range("E19").Value = "'2.55%"
$ws.Range("G19").Value = "'2"
$ws.Range("D20").Value = "'6.645"
$ws.Range("E20").Value = "'26.18%"
$ws.Range("G20").Value = "'2"
$ws.Range("D21").Value = "'0.1327"
$ws.Range("E21").Value = "'1.67%"
$ws.Range("G21").Value = "'2"
$ws.Range("D22").Value = "'0.2309"
$ws.Range("E22").Value = "'-10.85%"
$ws.Range("G22").Value = "'2"
$ws.Range("D23").Value = "'0.04414"
$ws.Range("E23").Value = "'1.57%"
$ws.Range("G23").Value = "'2"
$ws.Range("E24").Value = "'-1.65%"
$ws.Range("G24").Value = "'2"
$ws.Range("D25").Value = "'0.004361"
$ws.Range("E25").Value = "'-6.57%"
$ws.Range("G25").Value = "'2"
$ws.Range("D26").Value = "'0.0001139"
$ws.Range("E26").Value = "'-16.30%"
$ws.Range("G26").Value = "'2"
$ws.Range("D27").Value = "'0.0003986"
$ws.Range("E27").Value = "'-0.14%"
$ws.Range("G27").Value = "'2"
$ws.Range("G28").Value = "'2"
$ws.Range("G29").Value = "'2"
$ws.Range("G30").Value = "'2"
$ws.Range("G31").Value = "'2"
$ws.Range("G32").Value = "'2"
$ws.Range("G33").Value = "'2"
$ws.Range("G34").Value = "'2"
$ws.Range("G35").Value = "'2"
$ws.Range("G36").Value = "'2"
$ws.Range("G37").Value = "'2"
$ws.Range("G38").Value = "'2"
$ws.Range("D39").Value = "'0.02429"
$ws.Range("E39").Value = "'7.91%"
$ws.Range("G39").Value = "'2"
$ws.Range("D40").Value = "'0.05202"
$ws.Range("E40").Value = "'1.60%"
$ws.Range("G40").Value = "'2"
$ws.Range("D41").Value = "'0.007677"
$ws.Range("E41").Value = "'2.82%"
$ws.Range("G41").Value = "'2"
$ws.Range("D42").Value = "'0.1432"
$ws.Range("E42").Value = "'5.39%"
$ws.Range("G42").Value = "'2"
$ws.Range("D43").Value = "'0.008600"
$ws.Range("G43").Value = "'2"
$ws.Range("D44").Value = "'0.002108"
$ws.Range("E44").Value = "'5.34%"
$ws.Range("G44").Value = "'2"
$ws.Range("D45").Value = "'0.008981"
$ws.Range("E45").Value = "'4.04%"
$ws.Range("G45").Value = "'2"
$ws.Range("D46").Value = "'0.00006543"
$ws.Range("E46").Value = "'-1.28%"
$ws.Range("G46").Value = "'2"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.14%"
$ws.Range("G47").Value = "'2"
$ws.Range("D48").Value = "'0.002870"
$ws.Range("E48").Value = "'-11.51%"
$ws.Range("G48").Value = "'2"
$ws.Range("D49").Value = "'0.002489"
$ws.Range("E49").Value = "'148.81%"
$ws.Range("G49").Value = "'2"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").Value = "'-0.14%"
$ws.Range("G50").Value = "'2"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.14%"
$ws.Range("G51").Value = "'2"

# Remove the "number stored as text" formatting (quote-prefix style) that
# Excel applies when a numeric-looking literal is entered as text, so the
# cell style matches the original (unstyled) cells exactly.
$ws.Range("D2:G51").ClearFormats()
